# Update example with new speciation code.
# Renames the NH4 species bookkeeping to the new dissolved/adsorbed
# speciation split (TNH4_dis / TNH4_ads) across the model sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "species_in_model": rows 2 and 3 swap roles.
#     Row 2 becomes the N_org reactant, row 3 becomes the TNH4 product
#     with the updated species_eq list.
$ws1 = $wb.Worksheets.Item("species_in_model")

$ws1.Cells.Item(2,1).Value = "Rremin"
$ws1.Cells.Item(2,2).Value = "N_org"
$ws1.Cells.Item(2,3).Value = "-1"
$ws1.Cells.Item(2,4).Value = "0"
$ws1.Cells.Item(2,5).Value = "reactant"
$ws1.Cells.Item(2,6).Value = "N_org"
$ws1.Cells.Item(2,8).Value = "N_org"
$ws1.Cells.Item(2,9).Value = "solid"
$ws1.Cells.Item(2,10).Value = 1
$ws1.Cells.Item(2,11).Value = "solid"
$ws1.Cells.Item(2,12).Value = "solid"

$ws1.Cells.Item(3,1).Value = "Rremin"
$ws1.Cells.Item(3,2).Value = "TNH4"
$ws1.Cells.Item(3,3).Value = "1"
$ws1.Cells.Item(3,4).Value = "0"
$ws1.Cells.Item(3,5).Value = "product"
$ws1.Cells.Item(3,6).Value = "TNH4,NH4_ads,TNH4_ads,TNH4_ads_nsf,TNH4_dis"
$ws1.Cells.Item(3,8).Value = "TNH4"
$ws1.Cells.Item(3,9).Value = "dissolved_adsorbed_summed"
$ws1.Cells.Item(3,10).Value = 2
$ws1.Cells.Item(3,11).Value = "dissolved_adsorbed_summed"
$ws1.Cells.Item(3,12).Value = "solid"

# --- Sheet "reaction_dependency": rows 2 and 3 swap (TNH4 now depends
#     on N_org, in the second row; N_org on itself in the first).
$ws4 = $wb.Worksheets.Item("reaction_dependency")

$ws4.Cells.Item(2,1).Value = "N_org"
$ws4.Cells.Item(2,2).Value = "N_org"

$ws4.Cells.Item(3,1).Value = "TNH4"
$ws4.Cells.Item(3,2).Value = "N_org"

# --- Sheet "transport_parameters": NH4-prefixed identifiers become
#     TNH4_dis-prefixed ones; TNH4ID/N_orgID rows swap; final row (NH40,
#     adsorption) is removed because the bioirrigation/adsorption split
#     moved down and NH40 became TNH4_dis0 one row earlier; only 16 rows
#     remain.
$ws5 = $wb.Worksheets.Item("transport_parameters")

$ws5.Cells.Item(2,1).Value = "N_orgID"
$ws5.Cells.Item(2,2).Value = "index"

$ws5.Cells.Item(3,1).Value = "TNH4ID"
$ws5.Cells.Item(3,2).Value = "index"

$ws5.Cells.Item(5,1).Value = "AmTNH4_dis"
$ws5.Cells.Item(5,2).Value = "transport matrix"

$ws5.Cells.Item(6,1).Value = "AmTNH4_ads"
$ws5.Cells.Item(6,2).Value = "transport matrix"

$ws5.Cells.Item(10,1).Value = "BcAmTNH4_dis"
$ws5.Cells.Item(10,2).Value = "boundary condition"

$ws5.Cells.Item(11,1).Value = "BcCmTNH4_dis"
$ws5.Cells.Item(11,2).Value = "boundary condition"

$ws5.Cells.Item(12,1).Value = "BcAmTNH4_ads"
$ws5.Cells.Item(12,2).Value = "boundary condition"

$ws5.Cells.Item(13,1).Value = "BcCmTNH4_ads"
$ws5.Cells.Item(13,2).Value = "boundary condition"

$ws5.Cells.Item(15,1).Value = "TNH4_dis0"
$ws5.Cells.Item(15,2).Value = "bioirrigation"

$ws5.Cells.Item(16,1).Value = "dstopw"
$ws5.Cells.Item(16,2).Value = "adsorption"

$ws5.Rows.Item(17).Delete()

# --- Sheet "reaction_parameters": KNH4_ads/dstopw move from the
#     "adsorption" reaction type to a new "speciation" type, and their
#     comment formulas reference TNH4_dis instead of NH4.
$ws6 = $wb.Worksheets.Item("reaction_parameters")

$ws6.Cells.Item(2,1).Value = "KNH4_ads"
$ws6.Cells.Item(2,2).Value = "speciation"
$ws6.Cells.Item(2,3).Value = "TNH4_dis,NH4_ads"

$ws6.Cells.Item(3,1).Value = "dstopw"
$ws6.Cells.Item(3,2).Value = "speciation"
$ws6.Cells.Item(3,3).Value = "TNH4_dis"
